$d = $word.ActiveDocument

# The paragraph contains a run with text interspersed with manual line
# breaks (<w:br/>), including a trailing line break right before the
# paragraph mark. First remove that trailing line break entirely (it
# should not become a space), then replace all remaining manual line
# breaks with a single space, merging everything into one text run.

$p1 = $d.Paragraphs.Item(1)
$pr = $p1.Range
$n = $pr.Characters.Count

# Delete the last character of the paragraph's text (the manual line
# break that sits immediately before the paragraph mark).
$lastBreak = $d.Range($pr.Start + $n - 2, $pr.Start + $n - 1)
$lastBreak.Text = ""

# Replace every remaining manual line break (vertical tab, chr(11))
# with a single space across the whole document.
$d.Content.Find.Execute([char]11, $false, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null
